$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.037.00'
$ws.Range("E2").Value = '  +2.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.610.23'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.74'
$ws.Range("E5").Value = '  +10.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '567.06'
$ws.Range("E6").Value = '  -3.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.608.14'
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.621'
$ws.Range("E8").Value = '  +1.02%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("E11").Value = '  +15.13%  '
$ws.Range("E12").Value = '  +3.84%  '
$ws.Range("E13").Value = '  +11.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.02'
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.194.15'
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.614.77'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.87'
$ws.Range("E18").Value = '  +3.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.965.75'
$ws.Range("E19").Value = '  +2.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.38'
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.07'
$ws.Range("E21").Value = '  +2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '402.44'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.08'
$ws.Range("E23").Value = '  +17.21%  '
$ws.Range("E24").Value = '  -4.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.38'
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.95'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.63'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.95'
$ws.Range("E28").Value = '  +11.87%  '
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("E30").Value = '  +20.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.44'
$ws.Range("E31").Value = '  +5.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.66'
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '676.19'
$ws.Range("E33").Value = '  +8.48%  '
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '63.90'
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.24'
$ws.Range("E37").Value = '  +2.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.422'
$ws.Range("E38").Value = '  +6.30%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0771'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.272.21'
$ws.Range("E41").Value = '  +8.47%  '
$ws.Range("B42").Value = 'ThetaToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("E42").Value = '  +14.20%  '
$ws.Range("E43").Value = '  +4.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +10.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.07'
$ws.Range("E45").Value = '  +32.24%  '
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0419'
$ws.Range("E47").Value = '  +2.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.74'
$ws.Range("E48").Value = '  +10.93%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.85'
$ws.Range("E49").Value = '  +3.48%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.131'
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.05'
$ws.Range("E51").Value = '  +1.48%  '
